# "some changes on Training plan"
#
# 1) Profile / summary paragraph: append a new closing sentence about
#    looking for an overseas job.
# 2) "Training plan" bullet about the entertainment project: drop the
#    trailing whitespace runs left at the end of the paragraph.
# 3) NLU project bullet: reword "end users can control the internet TV
#    by using their voices" -> "end users can use their voices to
#    interact with".
# 4) Footer: the cached PAGE field result goes from 1 to 2.

$d = $word.ActiveDocument

# --- 1) Append sentence after the "...who speak English or French in
#        term of technique." paragraph -------------------------------
$d.Content.Find.Execute(
    "who speak English or French in term of technique.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "who speak English or French in term of technique. Actually, I" + [char]0x2019 + "m looking for an overseas job.",
    2) | Out-Null

# --- 2) Remove trailing whitespace after "entertainment project" -----
$d.Content.Find.Execute(
    "entertainment project   ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "entertainment project",
    2) | Out-Null

# --- 3) Reword the NLU voice-control sentence -------------------------
$d.Content.Find.Execute(
    "end users can control the internet TV by using their voices",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "end users can use their voices to interact with",
    2) | Out-Null

# --- 4) Footer page-number cached field text: 1 -> 2 -------------------
$sec = $d.Sections(1)
$footer = $sec.Footers(1)
$footer.Range.Find.Execute(
    "1", $true, $false, $false, $false, $false, $true, 1, $false,
    "2",
    2) | Out-Null
